# Fixed update to excel issue
#
# 1. Rename the generic "Requested quantity" header on the existing sheets
#    to metric-specific names.
# 2. Add a new "PO Forecast" sheet (placed right after "Monthly Trend")
#    holding the Prophet-style forecast output (ds / PO_Forecast /
#    yhat_lower / yhat_upper).

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new worksheet right after "Monthly Trend" ---------------------
$wsForecast = $wb.Worksheets.Add($null, $wsMonthly)
$wsForecast.Name = "PO Forecast"

# Re-use the existing header / date-column formatting so the new sheet
# matches the look of the other two (bold+bordered header row, date format
# on column A) before writing any values into it.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)   # xlPasteFormats

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A25").PasteSpecial(-4122)  # xlPasteFormats

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$data = @(
    @(44934.99999999999, 28, 17.83002544832913, 38.89439161451983),
    @(44941.99999999999, 28, 17.65333488876505, 37.88570649059221),
    @(44948.99999999999, 28, 17.10354212459917, 38.54420154734277),
    @(44955.99999999999, 28, 17.87589749292921, 38.16519901077096),
    @(44976.99999999999, 27, 16.27888044450215, 37.48528377966313),
    @(44983.99999999999, 27, 16.46997981463884, 36.94529469964278),
    @(44990.99999999999, 27, 15.51201814943122, 36.98513757749134),
    @(45004.99999999999, 26, 14.60345644479395, 37.40011466642191),
    @(45032.99999999999, 25, 15.09817765358964, 35.90300806255426),
    @(45039.99999999999, 25, 14.91618472336858, 35.93812003390864),
    @(45074.99999999999, 24, 13.92125376337776, 34.25981744354475),
    @(45095.99999999999, 23, 12.89846105359529, 33.58433324070339),
    @(45116.99999999999, 23, 11.62792834584737, 32.59053695175309),
    @(45151.99999999999, 22, 10.33702752709638, 30.91474309645853),
    @(45165.99999999999, 21, 9.700498035197981, 31.20876473485083),
    @(45221.99999999999, 19, 8.244719761223452, 28.97888159627795),
    @(45228.99999999999, 19, 8.544255738465077, 28.80611129814302),
    @(45235.99999999999, 19, 8.944170345706688, 29.81560253971754),
    @(45242.99999999999, 19, 9.003801585634584, 28.98665680038836),
    @(45249.99999999999, 19, 8.274836890280646, 29.4369139997879),
    @(45256.99999999999, 18, 8.774623421383772, 28.01806996285346),
    @(45263.99999999999, 18, 8.008547296245041, 28.01798157467684),
    @(45270.99999999999, 18, 7.721521685368502, 28.61299279560692),
    @(45277.99999999999, 18, 7.785041749498096, 28.7419585987625)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
